# feat: add 2022-Q1 data
#
# The workbook currently has sheets: 2020-Q4, 2021-Q1, 2021-Q3, 2021-Q4, 总计
# (总计 = "Total" summary sheet).
#
# After this edit it must have: 2020-Q4, 2021-Q1, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#   - "2022-Q1" is a brand-new per-quarter fund-holding sheet, inserted right
#     before "总计" (reusing the old 总计 sheet's position/sheetId).
#   - "总计" becomes a brand-new sheet (after 2022-Q1) holding the same kind
#     of summary table as before, plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# A template sheet that already has the exact "fund holding" layout/styling
# we need to reproduce for the new "2022-Q1" sheet (header row bold+bordered,
# index column A bold+bordered).
$template = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet.
# This keeps it positioned right after "2021-Q4" (where 总计 used to be).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Copy header-row + index-column formatting from the template sheet so the
# new sheet matches the look of the other per-quarter sheets.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2:A9").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund-holding rows (A is a 0-based row index, B..G are text, H is numeric)
$q1Data = @(
    @("0", "003984", "嘉实新能源新材料股票A", "77.01", "88.90", "4.42", "3.4038", 9),
    @("1", "110025", "易方达资源行业混合",     "21.23", "92.20", "5.18", "1.0997", 5),
    @("2", "003985", "嘉实新能源新材料股票C", "15.23", "88.90", "4.42", "0.6732", 9),
    @("3", "002657", "招商安裕灵活配置混合A", "32.90", "23.56", "1.11", "0.3652", 7),
    @("4", "005434", "鹏华睿投灵活配置混合",   "3.41",  "82.48", "1.80", "0.0614", 9),
    @("5", "002658", "招商安裕灵活配置混合C", "4.75",  "23.56", "1.11", "0.0527", 7),
    @("6", "002149", "嘉实新优选灵活配置混合", "0.22",  "93.76", "5.69", "0.0125", 7),
    @("7", "015206", "招商安裕灵活配置混合D", "0.88",  "23.56", "1.11", "0.0098", 7)
)

for ($i = 0; $i -lt $q1Data.Count; $i++) {
    $r = $i + 2
    $row = $q1Data[$i]
    $q1.Cells.Item($r, 1).Value = [int]$row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# Step 2: add a brand-new "总计" sheet right after "2022-Q1", rebuilding
# the summary table with a new leading row for 2022-Q1.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# Copy header-row + index-column formatting from the template sheet again.
$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$template.Range("A2:A6").Copy()
$total.Range("A2:A6").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @("0", "2022-Q1", 8,  5.68),
    @("1", "2021-Q4", 20, 8.029999999999999),
    @("2", "2021-Q3", 9,  1.83),
    @("3", "2021-Q1", 1,  0.03),
    @("4", "2020-Q4", 8,  4.55)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $r = $i + 2
    $row = $totalData[$i]
    $total.Cells.Item($r, 1).Value = [int]$row[0]
    $total.Cells.Item($r, 2).Value = "'" + $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

# Restore the originally active sheet/selection (adding sheets shifts focus
# onto the newly created one).
$wb.Worksheets.Item("2020-Q4").Activate()


